$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D (Price) and E (Volume(1h)) are treated as text so that
# values such as "1.000" or "30.135.71" are preserved exactly as strings.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '30.135.71'
$ws.Range("E2").Value = '  +4.58%  '

$ws.Range("D3").Value = '1.910.45'
$ws.Range("E3").Value = '  +5.45%  '

$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '251.49'
$ws.Range("E5").Value = '  +0.82%  '

$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.06%  '

$ws.Range("D7").Value = '0.5091'
$ws.Range("E7").Value = '  +2.90%  '

$ws.Range("D8").Value = '44.66'
$ws.Range("E8").Value = '  +3.21%  '

$ws.Range("D9").Value = '0.2955'
$ws.Range("E9").Value = '  +6.22%  '

$ws.Range("D10").Value = '0.06775'
$ws.Range("E10").Value = '  +5.91%  '

$ws.Range("D11").Value = '1.912.79'
$ws.Range("E11").Value = '  +5.60%  '

$ws.Range("D12").Value = '17.24'
$ws.Range("E12").Value = '  +2.96%  '

$ws.Range("D13").Value = '0.07363'
$ws.Range("E13").Value = '  +3.25%  '

$ws.Range("D14").Value = '0.6900'
$ws.Range("E14").Value = '  +6.42%  '

$ws.Range("D15").Value = '86.56'
$ws.Range("E15").Value = '  +3.01%  '

$ws.Range("D16").Value = '4.876'
$ws.Range("E16").Value = '  +3.90%  '

$ws.Range("D17").Value = '30.126.04'
$ws.Range("E17").Value = '  +4.65%  '

$ws.Range("D18").Value = '0.000008117'
$ws.Range("E18").Value = '  +9.72%  '

$ws.Range("E19").Value = '  +0.11%  '

$ws.Range("D20").Value = '12.97'
$ws.Range("E20").Value = '  +6.07%  '

$ws.Range("D21").Value = '2.158.28'
$ws.Range("E21").Value = '  +5.20%  '

$ws.Range("D22").Value = '0.9993'
$ws.Range("E22").Value = '  -0.06%  '

$ws.Range("D23").Value = '4.821'
$ws.Range("E23").Value = '  +4.92%  '

$ws.Range("D24").Value = '5.729'
$ws.Range("E24").Value = '  +7.24%  '

$ws.Range("D25").Value = '9.140'
$ws.Range("E25").Value = '  +2.85%  '

$ws.Range("D26").Value = '146.77'
$ws.Range("E26").Value = '  +2.63%  '

$ws.Range("D27").Value = '135.57'
$ws.Range("E27").Value = '  +1.79%  '

$ws.Range("D28").Value = '17.04'
$ws.Range("E28").Value = '  +2.63%  '

$ws.Range("D29").Value = '1.992'
$ws.Range("E29").Value = '  +5.54%  '

$ws.Range("D30").Value = '1.390'
$ws.Range("E30").Value = '  -0.68%  '

$ws.Range("E31").Value = '  +1.47%  '

$ws.Range("D32").Value = '0.08794'
$ws.Range("E32").Value = '  +5.26%  '

$ws.Range("D33").Value = '4.005'
$ws.Range("E33").Value = '  +4.10%  '

$ws.Range("D34").Value = '0.05061'
$ws.Range("E34").Value = '  +2.66%  '

$ws.Range("D35").Value = '1.145'
$ws.Range("E35").Value = '  +5.17%  '

$ws.Range("D36").Value = '0.7131'
$ws.Range("E36").Value = '  +5.34%  '

$ws.Range("D37").Value = '2.691'
$ws.Range("E37").Value = '  -0.71%  '

$ws.Range("D38").Value = '2.811'
$ws.Range("E38").Value = '  +1.73%  '

$ws.Range("D39").Value = '2.269'
$ws.Range("E39").Value = '  +0.10%  '

$ws.Range("E40").Value = '  +1.59%  '

$ws.Range("E41").Value = '  +6.51%  '

$ws.Range("D42").Value = '6.138'
$ws.Range("E42").Value = '  +1.38%  '

$ws.Range("D43").Value = '0.4285'
$ws.Range("E43").Value = '  +4.70%  '

$ws.Range("D44").Value = '104.95'
$ws.Range("E44").Value = '  +4.27%  '

$ws.Range("D45").Value = '0.9989'
$ws.Range("E45").Value = '  -0.08%  '

$ws.Range("D46").Value = '7.595'
$ws.Range("E46").Value = '  +5.51%  '

$ws.Range("D47").Value = '0.1278'
$ws.Range("E47").Value = '  +4.59%  '

$ws.Range("D48").Value = '0.05734'
$ws.Range("E48").Value = '  +3.95%  '

$ws.Range("D49").Value = '33.06'
$ws.Range("E49").Value = '  +4.54%  '

$ws.Range("D50").Value = '8.405'
$ws.Range("E50").Value = '  +3.02%  '

$ws.Range("D51").Value = '0.3794'
$ws.Range("E51").Value = '  +4.70%  '
